$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.738.81'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').Value = '1.637.28'
$ws.Range('E3').Value = '  +2.00%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'212.70"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('E6').Value = '  +1.84%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('D9').Value = "'0.0624"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.96%  '
$ws.Range('E10').Value = '  +5.01%  '
$ws.Range('D11').Value = "'0.0836"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.75%  '
$ws.Range('D12').Value = '1.867.65'
$ws.Range('E12').Value = '  +2.12%  '
$ws.Range('D13').Value = '1.643.85'
$ws.Range('E13').Value = '  +2.49%  '
$ws.Range('D14').Value = "'4.07"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.63%  '
$ws.Range('E15').Value = '  +2.82%  '
$ws.Range('D16').Value = '26.737.33'
$ws.Range('E16').Value = '  +2.09%  '
$ws.Range('D17').Value = "'63.13"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.27%  '
$ws.Range('D18').Value = '0.0₃0741'
$ws.Range('E18').Value = '  +2.05%  '
$ws.Range('D19').Value = "'208.67"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +4.04%  '
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('E21').Value = '  +1.22%  '
$ws.Range('E22').Value = '  +1.43%  '
$ws.Range('E23').Value = '  +2.51%  '
$ws.Range('E24').Value = '  +3.68%  '
$ws.Range('D25').Value = "'146.54"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.60%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').Value = "'0.120"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').Value = "'6.76"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.23%  '
$ws.Range('E29').Value = '  +1.52%  '
$ws.Range('E30').Value = '  +5.55%  '
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('D32').Value = "'3.23"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.06%  '
$ws.Range('E33').Value = '  +1.56%  '
$ws.Range('E34').Value = '  +1.09%  '
$ws.Range('E35').Value = '  +0.78%  '
$ws.Range('D36').Value = '1.169.35'
$ws.Range('E36').Value = '  +0.33%  '
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('E38').Value = '  +3.35%  '
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('E40').Value = '  +1.97%  '
$ws.Range('E41').Value = '  +0.41%  '
$ws.Range('D42').Value = "'0.794"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.88%  '
$ws.Range('D43').Value = "'5.38"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.63%  '
$ws.Range('D44').Value = '1.776.50'
$ws.Range('E44').Value = '  +2.10%  '
$ws.Range('D45').Value = "'92.54"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.22%  '
$ws.Range('E46').Value = '  +2.90%  '
$ws.Range('E47').Value = '  +8.42%  '
$ws.Range('D48').Value = "'54.89"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.62%  '
$ws.Range('E49').Value = '  +1.48%  '
$ws.Range('D50').Value = "'0.410"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.61%  '
$ws.Range('D51').Value = "'7.54"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +4.54%  '
